$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain-looking numeric strings (e.g. "55.212.38",
# "61.80", "1.00") as TEXT in this workbook (no thousands separator - the
# dots are literal). Excel normally reinterprets a numeric-looking Value
# assignment as a real number (dropping trailing zeros, etc.), so force
# those specific cells to Text format first, preserving the original intent.
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.60"
$ws.Range("E42").Value = "  -7.89%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.625"
$ws.Range("E43").Value = "  -6.16%  "
$ws.Range("D2").Value = "55.189.35"
$ws.Range("E2").Value = "  -4.85%  "
$ws.Range("D3").Value = "2.878.58"
$ws.Range("E3").Value = "  -5.66%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "486.47"
$ws.Range("E5").Value = "  -6.22%  "
$ws.Range("D6").Value = "130.98"
$ws.Range("E6").Value = "  -7.74%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -6.81%  "
$ws.Range("E9").Value = "  -5.78%  "
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  -8.28%  "
$ws.Range("D11").Value = "0.344"
$ws.Range("E11").Value = "  -6.90%  "
$ws.Range("D12").Value = "3.373.49"
$ws.Range("E12").Value = "  -5.66%  "
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").Value = "25.25"
$ws.Range("E14").Value = "  -5.66%  "
$ws.Range("D15").Value = "0.0000155"
$ws.Range("E15").Value = "  -8.35%  "
$ws.Range("D16").Value = "55.184.33"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "5.93"
$ws.Range("E17").Value = "  -4.91%  "
$ws.Range("D18").Value = "2.880.37"
$ws.Range("E18").Value = "  -5.51%  "
$ws.Range("D19").Value = "12.28"
$ws.Range("E19").Value = "  -5.84%  "
$ws.Range("E20").Value = "  -6.81%  "
$ws.Range("D21").Value = "309.84"
$ws.Range("E21").Value = "  -8.45%  "
$ws.Range("D24").Value = "0.475"
$ws.Range("E24").Value = "  -5.43%  "
$ws.Range("D25").Value = "61.80"
$ws.Range("E25").Value = "  -4.99%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  -6.56%  "
$ws.Range("D28").Value = "0.0₃0829"
$ws.Range("E28").Value = "  -13.32%  "
$ws.Range("D29").Value = "6.26"
$ws.Range("E29").Value = "  -9.47%  "
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -9.17%  "
$ws.Range("E31").Value = "  -6.30%  "
$ws.Range("D32").Value = "19.48"
$ws.Range("E32").Value = "  -7.37%  "
$ws.Range("D33").Value = "1.10"
$ws.Range("E33").Value = "  -10.56%  "
$ws.Range("D34").Value = "148.36"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  -9.38%  "
$ws.Range("D36").Value = "5.49"
$ws.Range("E36").Value = "  -7.36%  "
$ws.Range("D37").Value = "24.16"
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("D38").Value = "1.16"
$ws.Range("E38").Value = "  -10.23%  "
$ws.Range("E39").Value = "  -7.08%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "35.96"
$ws.Range("E41").Value = "  -4.87%  "
$ws.Range("D44").Value = "2.075.19"
$ws.Range("E44").Value = "  -10.87%  "
$ws.Range("E45").Value = "  -9.83%  "
$ws.Range("D46").Value = "5.77"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("D47").Value = "0.897"
$ws.Range("E47").Value = "  -11.15%  "
$ws.Range("D48").Value = "0.0227"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("D49").Value = "18.41"
$ws.Range("E49").Value = "  -7.20%  "
$ws.Range("D50").Value = "0.0833"
$ws.Range("E50").Value = "  -7.38%  "
$ws.Range("E51").Value = "  -8.03%  "
